$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price column so numeric-looking strings are preserved exactly
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.559.14"
$ws.Range("E2").Value = "  +3.20%  "

$ws.Range("D3").Value = "3.067.11"
$ws.Range("E3").Value = "  +3.35%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "517.66"
$ws.Range("E5").Value = "  +3.48%  "

$ws.Range("D6").Value = "141.20"
$ws.Range("E6").Value = "  +3.57%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "0.433"
$ws.Range("E8").Value = "  +1.96%  "

$ws.Range("D9").Value = "7.22"
$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +1.83%  "

$ws.Range("D11").Value = "0.374"
$ws.Range("E11").Value = "  +3.59%  "

$ws.Range("D12").Value = "3.601.23"
$ws.Range("E12").Value = "  +3.47%  "

$ws.Range("D13").Value = "0.129"
$ws.Range("E13").Value = "  +3.21%  "

$ws.Range("D14").Value = "25.58"
$ws.Range("E14").Value = "  -0.88%  "

$ws.Range("D15").Value = "0.0000163"
$ws.Range("E15").Value = "  +2.36%  "

$ws.Range("D16").Value = "57.611.22"
$ws.Range("E16").Value = "  +3.23%  "

$ws.Range("D17").Value = "3.067.96"
$ws.Range("E17").Value = "  +3.52%  "

$ws.Range("D18").Value = "6.12"
$ws.Range("E18").Value = "  +2.65%  "

$ws.Range("D19").Value = "12.86"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").Value = "8.10"
$ws.Range("E20").Value = "  +2.04%  "

$ws.Range("D21").Value = "331.79"
$ws.Range("E21").Value = "  +1.47%  "

$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").Value = "0.498"
$ws.Range("E23").Value = "  +1.70%  "

$ws.Range("D24").Value = "65.68"
$ws.Range("E24").Value = "  +2.28%  "

$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +4.89%  "

$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "0.0₃0901"
$ws.Range("E27").Value = "  +1.91%  "

$ws.Range("D28").Value = "6.38"
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  +4.31%  "

$ws.Range("D30").Value = "1.82"
$ws.Range("E30").Value = "  +3.47%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +4.89%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "20.75"
$ws.Range("E32").Value = "  +3.32%  "

$ws.Range("D33").Value = "154.47"
$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  +0.74%  "

$ws.Range("D35").Value = "26.94"
$ws.Range("E35").Value = "  +7.84%  "

$ws.Range("D36").Value = "5.92"
$ws.Range("E36").Value = "  +4.70%  "

$ws.Range("D37").Value = "1.26"
$ws.Range("E37").Value = "  +2.96%  "

$ws.Range("D38").Value = "0.0672"
$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("D39").Value = "3.110.93"
$ws.Range("E39").Value = "  +3.65%  "

$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  +4.53%  "

$ws.Range("D41").Value = "36.80"
$ws.Range("E41").Value = "  +0.44%  "

$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").Value = "0.655"
$ws.Range("E43").Value = "  +1.36%  "

$ws.Range("D44").Value = "2.262.02"
$ws.Range("E44").Value = "  +5.39%  "

$ws.Range("D45").Value = "0.0259"
$ws.Range("E45").Value = "  +10.85%  "

$ws.Range("D46").Value = "20.78"
$ws.Range("E46").Value = "  +7.05%  "

$ws.Range("D47").Value = "1.36"
$ws.Range("E47").Value = "  +1.82%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "5.87"
$ws.Range("E48").Value = "  +1.35%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "0.926"
$ws.Range("E49").Value = "  +1.20%  "

$ws.Range("D50").Value = "263.17"
$ws.Range("E50").Value = "  +16.54%  "

$ws.Range("D51").Value = "0.709"
$ws.Range("E51").Value = "  +6.56%  "
